$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10
$ws.Range("A10").Value = "ARH0227"
$ws.Range("B10").Value = "Deweloper nieruchomości"
$ws.Range("C10").Value = 101.62
$ws.Range("D10").Value = 1
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "2025-01-07"
$ws.Range("E10").ClearFormats()

# Row 11
$ws.Range("A11").Value = "KRI1025"
$ws.Range("B11").Value = "Zarządzanie wierzytelnościami"
$ws.Range("C11").Value = 99.3
$ws.Range("D11").Value = 6
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "2025-01-07"
$ws.Range("E11").ClearFormats()
